# Add "sync force" & "sync rotation" interval columns to the configWorld sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("configWorld")

# Make configWorld the active sheet (matches the saved workbook's active tab).
$ws.Select() | Out-Null

# Insert two new blank columns before the existing "地图逻辑处理分格大小" column.
$ws.Columns("B:C").Insert()

# Row 1: Chinese display names.
$ws.Cells.Item(1, 2).Value = "同步移动间隔(帧)"
$ws.Cells.Item(1, 3).Value = "同步旋转间隔(帧)"

# Row 2: field/key names.
$ws.Cells.Item(2, 2).Value = "syncForceFrame"
$ws.Cells.Item(2, 3).Value = "syncRotationFrame"

# Row 3: default values.
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(3, 3).Value = 2

# Match the saved selection location.
$ws.Range("C10").Select() | Out-Null
